# Applies the "Add files via upload" edit to the Landscaping Data workbook:
#  - Update three Growth (column H) values in existing rows 294, 295, 302
#  - Append 7 new data rows (310-316) for 2025-06-23 (Excel serial 45831)
#  - Extend the F-column ABS(Low-High) formula down through the new rows
#  - Update the active selection to reflect the new bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Tweak existing Growth values
# ---------------------------------------------------------------------
$ws.Cells.Item(294, 8).Value = 0
$ws.Cells.Item(295, 8).Value = 1
$ws.Cells.Item(302, 8).Value = 0.5

# ---------------------------------------------------------------------
# 2) Append the new rows of data collected on 6/23/2025 (serial 45831)
#    Column order below (F is a formula column, filled in separately):
#    A Date, B Plant_Type, C Plant_Size, D Low, E High, G Rain, H Growth,
#    I Pruned, J Quadrant, K Shade, L UV, M Humidity, N Dew_Point,
#    O Pressure, P Wind_Gust, Q Cloud_Cover, R Visibility, S AQI, T Pollen
# ---------------------------------------------------------------------
$newRows = @(
    @(45831, "Flowering",    "Large",  72, 94, 0, 0,   "Yes", 2, "Bright",  9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11),
    @(45831, "Nonflowering", "Medium", 72, 94, 0, 0.1, "Yes", 3, "Neutral", 9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11),
    @(45831, "Nonflowering", "Small",  72, 94, 0, 0.1, "Yes", 3, "Bright",  9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11),
    @(45831, "Nonflowering", "Medium", 72, 94, 0, 0,   "Yes", 3, "Neutral", 9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11),
    @(45831, "Nonflowering", "Medium", 72, 94, 0, 0.2, "Yes", 3, "Bright",  9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11),
    @(45831, "Nonflowering", "Large",  72, 94, 0, 0.3, "Yes", 4, "Neutral", 9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11),
    @(45831, "Tree",         "Medium", 72, 94, 0, 0.5, "Yes", 1, "Bright",  9, 0.54, 73, 30.19, 10, 0.06, 9.9, 60, 11)
)

$startRow = 310
$r = $startRow
foreach ($row in $newRows) {
    # Seed the date cell by copying the style (number format) of the row
    # directly above it, then overwrite the value so no new numFmt is
    # created in styles.xml (keeps the same date style as the rest of col A).
    $ws.Cells.Item($r - 1, 1).Copy($ws.Cells.Item($r, 1))

    $ws.Cells.Item($r, 1).Value  = $row[0]      # A Date
    $ws.Cells.Item($r, 2).Value  = $row[1]      # B Plant_Type
    $ws.Cells.Item($r, 3).Value  = $row[2]      # C Plant_Size
    $ws.Cells.Item($r, 4).Value  = $row[3]      # D Low
    $ws.Cells.Item($r, 5).Value  = $row[4]      # E High
    $ws.Cells.Item($r, 6).Formula = "=ABS(D$r-E$r)"   # F Temp_Diff
    $ws.Cells.Item($r, 7).Value  = $row[5]      # G Rain
    $ws.Cells.Item($r, 8).Value  = $row[6]      # H Growth
    $ws.Cells.Item($r, 9).Value  = $row[7]      # I Pruned
    $ws.Cells.Item($r, 10).Value = $row[8]      # J Quadrant
    $ws.Cells.Item($r, 11).Value = $row[9]      # K Shade
    $ws.Cells.Item($r, 12).Value = $row[10]     # L UV
    $ws.Cells.Item($r, 13).Value = $row[11]     # M Humidity
    $ws.Cells.Item($r, 14).Value = $row[12]     # N Dew_Point
    $ws.Cells.Item($r, 15).Value = $row[13]     # O Pressure
    $ws.Cells.Item($r, 16).Value = $row[14]     # P Wind_Gust
    $ws.Cells.Item($r, 17).Value = $row[15]     # Q Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $row[16]     # R Visibility
    $ws.Cells.Item($r, 19).Value = $row[17]     # S AQI
    $ws.Cells.Item($r, 20).Value = $row[18]     # T Pollen

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Update the view/selection to match the edited file (scrolled to the
#    new bottom of the sheet, active cell just past the last data row).
# ---------------------------------------------------------------------
$ws.Range("I317").Select()
